# "Add files via upload" — refresh the submission-group rosters.
#
# Sheet "11" (group 1, previously empty) gets the four names that used to
# live in sheet "12"'s column C, but spelled out with this round's actual
# submitters. Sheet "12" keeps its own column A roster and loses the old
# column C roster entirely (those four names move to sheet "11").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("11")
$ws2 = $wb.Worksheets.Item("12")

# sheet "12": drop the old column-C roster (Liel Azulay / Talia Arbiv /
# Adi Zadicareo / Paz Elraz) — column A stays untouched.
$ws2.Range("C1:C4").ClearContents()

# sheet "11": populate column A with the new roster. A4 is written before
# A3 so the new shared-string entries land in the same order as the
# target workbook.
$ws1.Range("A1").Value = "batel elbaz"
$ws1.Range("A2").Value = "shahar gavriel"
$ws1.Range("A4").Value = "liad tzvaot"
$ws1.Range("A3").Value = "idan yontov"

# Restore the selections: sheet "12" lands on A4 (no longer the active
# tab), sheet "11" ends up active with E3 selected — select sheet "12"
# first so sheet "11" is the last (and therefore active) tab.
$ws2.Range("A4").Select()
$ws1.Range("E3").Select()
